# "no data for Disciplinary Actions 2003"
#
# The workbook had a standalone title cell in row 1 ("Criminal Offenses -
# On-campus Student Housing Facilities") with the real column headers
# living on row 2. The edit removes that title row entirely (so the header
# row becomes row 1 and every data row shifts up by one), and also
# normalizes the capitalization of several column header labels.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the title-only row. This shifts the header row (old row 2) up to
# row 1, and every data row (old rows 3..55) up by one as well.
$ws.Rows.Item(1).Delete()

# Re-label the header row (now row 1) with corrected capitalization.
$ws.Range("A1").Value = "Survey Year"
$ws.Range("B1").Value = "UnitID"
$ws.Range("C1").Value = "Institution Name"
$ws.Range("D1").Value = "Campus ID"
$ws.Range("E1").Value = "Campus Name"
$ws.Range("F1").Value = "Institution Size"
$ws.Range("G1").Value = "Murder/Non-Negligent Manslaughter"
$ws.Range("H1").Value = "Negligent Manslaughter"
$ws.Range("I1").Value = "Sex Offenses - Forcible"
$ws.Range("J1").Value = "Sex Offenses - Non-Forcible"
$ws.Range("K1").Value = "Robbery"
$ws.Range("L1").Value = "Aggravated Assault"
$ws.Range("M1").Value = "Burglary"
$ws.Range("N1").Value = "Motor Vehicle Theft"
$ws.Range("O1").Value = "Arson"
